$wb = $excel.ActiveWorkbook

$source = $wb.Sheets.Item("Turkey")
$source.Cells.Select()
$source.Copy([System.Reflection.Missing]::Value, $wb.Sheets.Item($wb.Sheets.Count))

$newSheet = $wb.Sheets.Item($wb.Sheets.Count)
$newSheet.Name = "Croatia"

$newSheet.Range("B2").Value = "Croatia Market"
$newSheet.Range("B4").Value = "NGC-3139/T2477"

$newSheet.Cells.UnMerge()

$newSheet.Range("B4").Select()
